$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 10; $row++) {
    $hCell = $ws.Cells.Item($row, 8)
    $iCell = $ws.Cells.Item($row, 9)

    # Decrement "PERIOD TO EXPIRE" by one day.
    $hCell.Value2 = $hCell.Value2 - 1

    # Update "LAST UPDATE" text to the new date, keeping it as plain text
    # (not an Excel date serial) and preserving the cell's original style.
    $iCell.NumberFormat = "@"
    $iCell.Value2 = "04-Nov-2025"
    $hCell.Copy()
    $iCell.PasteSpecial(-4122) # xlPasteFormats
}

$excel.CutCopyMode = $false
